$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 62763.062
$ws.Range("I103").Value = 77160.53999999999
$ws.Range("J103").Value = 374
$ws.Range("K103").Value = 231481.62
$ws.Range("L103").Value = 1122
$ws.Range("M103").Value = -230895.62
$ws.Range("N103").Value = -2294

$ws.Range("H138").Value = 2727.717
$ws.Range("I138").Value = 4176.3335
$ws.Range("J138").Value = 2431.4092
$ws.Range("K138").Value = 12529.0005
$ws.Range("L138").Value = 7294.2276
$ws.Range("M138").Value = -7389.000499999998
$ws.Range("N138").Value = -17574.2276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 14242.6
$ws.Range("I15").Value = 300
$ws.Range("K15").Value = 300
$ws.Range("M15").Value = 50

$ws.Range("H61").Value = 10755156
$ws.Range("I61").Value = 27779278
$ws.Range("J61").Value = 3078.5264
$ws.Range("K61").Value = 27779278
$ws.Range("L61").Value = 3078.5264
$ws.Range("M61").Value = -27779066
$ws.Range("N61").Value = -3502.5264

$ws.Range("H74").Value = 1108.3636
$ws.Range("I74").Value = 770.125
$ws.Range("J74").Value = 1301.6428
$ws.Range("K74").Value = 770.125
$ws.Range("L74").Value = 1301.6428
$ws.Range("M74").Value = 103.875
$ws.Range("N74").Value = -3049.6428

$ws.Range("H77").Value = 1108.3636
$ws.Range("I77").Value = 770.125
$ws.Range("J77").Value = 1301.6428
$ws.Range("K77").Value = 3850.625
$ws.Range("L77").Value = 6508.214
$ws.Range("M77").Value = 517.375
$ws.Range("N77").Value = -15244.214

$ws.Range("H88").Value = 3112
$ws.Range("I88").Value = 3168
$ws.Range("K88").Value = 3168
$ws.Range("M88").Value = -2762

$ws.Range("H91").Value = 3112
$ws.Range("I91").Value = 3168
$ws.Range("K91").Value = 3168
$ws.Range("M91").Value = -1764

$ws.Range("H110").Value = 1854.4783
$ws.Range("I110").Value = 2047.4
$ws.Range("J110").Value = 568.3333
$ws.Range("K110").Value = 2047.4
$ws.Range("L110").Value = 568.3333
$ws.Range("M110").Value = -2.400000000000091
$ws.Range("N110").Value = -4658.3333

$ws.Range("H136").Value = 10755156
$ws.Range("I136").Value = 27779278
$ws.Range("J136").Value = 3078.5264
$ws.Range("K136").Value = 83337834
$ws.Range("L136").Value = 9235.5792
$ws.Range("M136").Value = -83335284
$ws.Range("N136").Value = -14335.5792

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 44999
$ws.Range("J58").Value = 44999
$ws.Range("L58").Value = 44999
$ws.Range("N58").Value = -45587

$ws.Range("H60").Value = 30333
$ws.Range("J60").Value = 30333
$ws.Range("L60").Value = 30333
$ws.Range("N60").Value = -31531

$ws.Range("H105").Value = 1962.64
$ws.Range("I105").Value = 1492.1177
$ws.Range("J105").Value = 2962.5
$ws.Range("K105").Value = 1492.1177
$ws.Range("L105").Value = 2962.5
$ws.Range("M105").Value = 254.8823
$ws.Range("N105").Value = -6456.5

$ws.Range("H134").Value = 3147.5334
$ws.Range("I134").Value = 3199.2942
$ws.Range("K134").Value = 9597.882599999999
$ws.Range("M134").Value = -7062.882599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2841.2856
$ws.Range("I58").Value = 2727.8
$ws.Range("K58").Value = 2727.8
$ws.Range("M58").Value = -2524.8

$ws.Range("H136").Value = 2841.2856
$ws.Range("I136").Value = 2727.8
$ws.Range("K136").Value = 8183.400000000001
$ws.Range("M136").Value = -5633.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1359
$ws.Range("I68").Value = 810.7241
$ws.Range("J68").Value = 1582.9436
$ws.Range("K68").Value = 2432.1723
$ws.Range("L68").Value = 4748.8308
$ws.Range("M68").Value = -1621.1723
$ws.Range("N68").Value = -6370.8308

$ws.Range("H71").Value = 1359
$ws.Range("I71").Value = 810.7241
$ws.Range("J71").Value = 1582.9436
$ws.Range("K71").Value = 7296.516900000001
$ws.Range("L71").Value = 14246.4924
$ws.Range("M71").Value = -3240.516900000001
$ws.Range("N71").Value = -22358.4924

$ws.Range("H107").Value = 1514.5
$ws.Range("I107").Value = 317.93332
$ws.Range("J107").Value = 2330.3408
$ws.Range("K107").Value = 953.7999599999999
$ws.Range("L107").Value = 6991.0224
$ws.Range("M107").Value = 966.2000400000001
$ws.Range("N107").Value = -10831.0224

$ws.Range("H131").Value = 1033.5
$ws.Range("J131").Value = 1170
$ws.Range("L131").Value = 3510
$ws.Range("N131").Value = -13590

$ws.Range("H134").Value = 2870.3635
$ws.Range("I134").Value = 2362.45
$ws.Range("J134").Value = 7949.5
$ws.Range("K134").Value = 7087.349999999999
$ws.Range("L134").Value = 23848.5
$ws.Range("M134").Value = -2017.349999999999
$ws.Range("N134").Value = -33988.5

$ws.Range("H140").Value = 1974.96
$ws.Range("J140").Value = 2892.9
$ws.Range("L140").Value = 8678.700000000001
$ws.Range("N140").Value = -19038.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -515
$ws.Range("N36").Value = $null

$ws.Range("H42").Value = 66145
$ws.Range("J42").Value = 66145
$ws.Range("L42").Value = 66145
$ws.Range("N42").Value = -67115

$ws.Range("H46").Value = 4285.6665
$ws.Range("J46").Value = 4285.6665
$ws.Range("L46").Value = 4285.6665
$ws.Range("N46").Value = -4597.6665

$ws.Range("H70").Value = 5957.269
$ws.Range("I70").Value = 5875.294
$ws.Range("J70").Value = 6112.1113
$ws.Range("K70").Value = 5875.294
$ws.Range("L70").Value = 6112.1113
$ws.Range("M70").Value = -5605.294
$ws.Range("N70").Value = -6652.1113

$ws.Range("H73").Value = 5957.269
$ws.Range("I73").Value = 5875.294
$ws.Range("J73").Value = 6112.1113
$ws.Range("K73").Value = 5875.294
$ws.Range("L73").Value = 6112.1113
$ws.Range("M73").Value = -4939.294
$ws.Range("N73").Value = -7984.1113

$ws.Range("H102").Value = 2415
$ws.Range("I102").Value = 2375
$ws.Range("J102").Value = 2495
$ws.Range("K102").Value = 2375
$ws.Range("L102").Value = 2495
$ws.Range("M102").Value = -753
$ws.Range("N102").Value = -5739

$ws.Range("H115").Value = 66145
$ws.Range("J115").Value = 66145
$ws.Range("L115").Value = 66145
$ws.Range("N115").Value = -68495

$ws.Range("H122").Value = 5357.2856
$ws.Range("I122").Value = 6004.6665
$ws.Range("J122").Value = 5180.727
$ws.Range("K122").Value = 18013.9995
$ws.Range("L122").Value = 15542.181
$ws.Range("M122").Value = -15563.9995
$ws.Range("N122").Value = -20442.181

$ws.Range("H132").Value = 2023.3462
$ws.Range("I132").Value = 1785.75
$ws.Range("J132").Value = 2227
$ws.Range("K132").Value = 5357.25
$ws.Range("L132").Value = 6681
$ws.Range("M132").Value = -2827.25
$ws.Range("N132").Value = -11741

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 200002260
$ws.Range("I7").Value = 250001570
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 250001570
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -250001458
$ws.Range("N7").Value = -5224

$ws.Range("H42").Value = 4000
$ws.Range("J42").Value = 4000
$ws.Range("L42").Value = 4000
$ws.Range("N42").Value = -5126

$ws.Range("H49").Value = 4000
$ws.Range("J49").Value = 4000
$ws.Range("L49").Value = 4000
$ws.Range("N49").Value = -4294

$ws.Range("H126").Value = 200002260
$ws.Range("I126").Value = 250001570
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 750004710
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -750002240
$ws.Range("N126").Value = -19940

$ws.Range("H132").Value = 2547.22
$ws.Range("I132").Value = 2111.2058
$ws.Range("J132").Value = 3473.75
$ws.Range("K132").Value = 6333.617400000001
$ws.Range("L132").Value = 10421.25
$ws.Range("M132").Value = -3803.617400000001
$ws.Range("N132").Value = -15481.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = $null

$ws.Range("H136").Value = 2240.35
$ws.Range("I136").Value = 2137.068
$ws.Range("J136").Value = 2524.375
$ws.Range("K136").Value = 6411.204000000001
$ws.Range("L136").Value = 7573.125
$ws.Range("M136").Value = -3861.204000000001
$ws.Range("N136").Value = -12673.125
